$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: unmerge every merged range whose shape changes, BEFORE touching values/styles ---
$ws.Range("B5:B8").UnMerge()
$ws.Range("A5:A8").UnMerge()
$ws.Range("B11").UnMerge()
$ws.Range("A11").UnMerge()
$ws.Range("B12:B13").UnMerge()
$ws.Range("A12:A13").UnMerge()
$ws.Range("B14:B15").UnMerge()
$ws.Range("A14:A15").UnMerge()

# --- Step 2: header / info block edits ---
$ws.Range("A1").Value = "firma54"
$ws.Range("B2").Value = "19.3.2020"
$ws.Range("E2").Value = "ec22/2020"
$ws.Range("B3").Value = "1.1.1"
$ws.Range("E3").Value = 0

# --- Step 3: data rows 5-12 ---
# Row 5
$ws.Range("A5").Value = "M1"
$ws.Range("B5").Value = 78

# Row 6
$ws.Range("C6").Value = 44
$ws.Range("D6").Value = "czarny"

# Row 7
$ws.Range("C7").Value = 22
$ws.Range("D7").Value = "bialy"

# Row 8 - becomes a new data "group head" row (style like row 9's A/B cells)
$ws.Range("A8").Style = $ws.Range("A9").Style
$ws.Range("B8").Style = $ws.Range("B9").Style
$ws.Range("A8").Value = "Statyw metalowy"
$ws.Range("B8").Value = 33
$ws.Range("C8").Value = 33
$ws.Range("D8").Value = "90"

# Row 9
$ws.Range("A9").Value = "Akcesoria"
$ws.Range("B9").Value = 45
$ws.Range("C9").Value = 12
$ws.Range("D9").Value = "trzpień"

# Row 10
$ws.Range("C10").Value = 33
$ws.Range("D10").Value = "trzpień"

# Row 11
$ws.Range("A11").Value = "Statyw drewniany"
$ws.Range("B11").Value = 567
$ws.Range("C11").Value = 12
$ws.Range("D11").Value = "biały"

# Row 12 - loses its own group-head data, becomes a blank continuation row (style like row 7)
$ws.Range("A12").Style = $ws.Range("A7").Style
$ws.Range("B12").Style = $ws.Range("B7").Style
$ws.Range("A12").ClearContents()
$ws.Range("B12").ClearContents()
$ws.Range("C12").Value = 555
$ws.Range("D12").Value = "czarny"

# --- Step 4: rows 13-15 become fully empty rows ---
$ws.Range("A13:E13").ClearContents()
$ws.Range("A14:E14").ClearContents()
$ws.Range("A15:E15").ClearContents()
$ws.Range("A13:E13").Style = $ws.Range("A16:E16").Style
$ws.Range("A14:E14").Style = $ws.Range("A16:E16").Style
$ws.Range("A15:E15").Style = $ws.Range("A16:E16").Style

# --- Step 5: re-create the merges matching the new layout ---
$ws.Range("B5:B7").Merge()
$ws.Range("A5:A7").Merge()
$ws.Range("B8").Merge()
$ws.Range("A8").Merge()
$ws.Range("B11:B12").Merge()
$ws.Range("A11:A12").Merge()
